$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3008579558616162
$ws.Range("C2").Value = 0.8765574948681092
$ws.Range("D2").Value = 1.952297331720331
$ws.Range("E2").Value = 1.510493160416715
$ws.Range("F2").Value = 0.7007203280930838
$ws.Range("G2").Value = 0.9808517754788966
$ws.Range("H2").Value = 1.019718840708109
$ws.Range("B3").Value = 0.575699539006493
$ws.Range("C3").Value = 1.651439375858714
$ws.Range("D3").Value = 1.209635204555099
$ws.Range("E3").Value = 0.3998623722314676
$ws.Range("F3").Value = 0.6799938196172803
$ws.Range("G3").Value = 0.7188608848464924
$ws.Range("B4").Value = 1.075739836852222
$ws.Range("C4").Value = 0.6339356655486057
$ws.Range("D4").Value = -0.1758371667750253
$ws.Range("E4").Value = 0.1042942806107874
$ws.Range("F4").Value = 0.1431613458399994
$ws.Range("G4").Value = 0.03444913272475958
$ws.Range("H4").Value = 0.1946995151992468
$ws.Range("I4").Value = 0.1026476495847817
$ws.Range("J4").Value = -0.1135594886937667
$ws.Range("B5").Value = -0.4418041713036158
$ws.Range("C5").Value = -1.251577003627247
$ws.Range("D5").Value = -0.9714455562414341
$ws.Range("E5").Value = -0.9325784910122221
$ws.Range("F5").Value = -1.041290704127462
$ws.Range("G5").Value = -0.8810403216529747
$ws.Range("H5").Value = -0.9730921872674398
$ws.Range("I5").Value = -1.189299325545988
$ws.Range("B6").Value = -0.8097728323236311
$ws.Range("C6").Value = -0.5296413849378183
$ws.Range("D6").Value = -0.4907743197086063
$ws.Range("E6").Value = -0.5994865328238461
$ws.Range("F6").Value = -0.4392361503493589
$ws.Range("G6").Value = -0.531288015963824
$ws.Range("H6").Value = -0.7474951542423725
$ws.Range("B7").Value = 0.2801314473858127
$ws.Range("C7").Value = 0.3189985126150248
$ws.Range("D7").Value = 0.2102862994997849
$ws.Range("E7").Value = 0.3705366819742721
$ws.Range("F7").Value = 0.278484816359807
$ws.Range("G7").Value = 0.06227767808125861
$ws.Range("B8").Value = 0.03886706522921202
$ws.Range("C8").Value = -0.0698451478860278
$ws.Range("D8").Value = 0.09040523458845939
$ws.Range("E8").Value = -0.0016466310260057
$ws.Range("F8").Value = -0.2178537693045541
$ws.Range("G8").Value = -0.5353744557447726
$ws.Range("H8").Value = -0.2431494525332169
$ws.Range("I8").Value = -0.2492947008739179
$ws.Range("B9").Value = -0.1087122131152398
$ws.Range("C9").Value = 0.05153816935924738
$ws.Range("D9").Value = -0.04051369625521772
$ws.Range("E9").Value = -0.2567208345337662
$ws.Range("F9").Value = -0.5742415209739846
$ws.Range("G9").Value = -0.2820165177624289
$ws.Range("H9").Value = -0.2881617661031299
$ws.Range("B10").Value = 0.1602503824744872
$ws.Range("C10").Value = 0.0681985168600221
$ws.Range("D10").Value = -0.1480086214185263
$ws.Range("E10").Value = -0.4655293078587448
$ws.Range("F10").Value = -0.1733043046471891
$ws.Range("G10").Value = -0.1794495529878901
$ws.Range("B11").Value = -0.0920518656144651
$ws.Range("C11").Value = -0.3082590038930135
$ws.Range("D11").Value = -0.625779690333232
$ws.Range("E11").Value = -0.3335546871216763
$ws.Range("F11").Value = -0.3396999354623773
$ws.Range("B12").Value = -0.2162071382785484
$ws.Range("C12").Value = -0.5337278247187669
$ws.Range("D12").Value = -0.2415028215072112
$ws.Range("E12").Value = -0.2476480698479122
$ws.Range("B13").Value = -0.3175206864402185
$ws.Range("C13").Value = -0.0252956832286628
$ws.Range("D13").Value = -0.03144093156936381
$ws.Range("B14").Value = 0.2922250032115557
$ws.Range("C14").Value = 0.2860797548708547
$ws.Range("B15").Value = -0.006145248340701015
